$wb = $excel.ActiveWorkbook

$sheetNames = @("HVAC", "ARCHITECTURE")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 19; $row++) {
        $ws.Cells.Item($row, 2).Value = "0"
    }
    $ws.Range("B19").Select()
}
